$d = $word.ActiveDocument

$rng = $d.Content
$ok = $rng.Find.Execute("Q1: What are the applications that are getting impacted with that?")
if (-not $ok) { Write-Output "NOT FOUND #1: Q1: What are the applications that are getting impacted with that?" } else { $rng.Text = "Q1: What is the issue at the high critical site in the US?" }

$rng = $d.Content
$ok = $rng.Find.Execute("A1: All the applications for that one site are impacted.")
if (-not $ok) { Write-Output "NOT FOUND #2: A1: All the applications for that one site are impacted." } else { $rng.Text = "A1: It's a factory site, and both AB are not reachable. Praveen confirmed a power issue at the site. On-site personnel are checking the circuit breaker related to building electricity." }

$rng = $d.Content
$ok = $rng.Find.Execute("Q2: What is the issue with the high critical site located in the US?")
if (-not $ok) { Write-Output "NOT FOUND #3: Q2: What is the issue with the high critical site located in the US?" } else { $rng.Text = "Q2: What is the current status of the network issue after the power failure?" }

$rng = $d.Content
$ok = $rng.Find.Execute("A2: The site is a factory site, and both of the AB are not reachable. There is a power issue confirmed by Praveen, which has caused the site to become unreachable.")
if (-not $ok) { Write-Output "NOT FOUND #4: A2: The site is a factory site, and both of the AB are not reachable. There is a power issue confirmed by Praveen, which has caused the site to become unreachable." } else { $rng.Text = "A2: The main power in the building is up, but there is one room where the power has failed, affecting network equipment. The electricians are working to resolve this, which is why the whole network is down." }

$rng = $d.Content
$ok = $rng.Find.Execute("Q3: Is the site having just one power source, or are the devices that went down connected to one power source?")
if (-not $ok) { Write-Output "NOT FOUND #5: Q3: Is the site having just one power source, or are the devices that went down connected to one power source?" } else { $rng.Text = "Q3: What steps are being taken to resolve the network issue?" }

$rng = $d.Content
$ok = $rng.Find.Execute("A3: We'll have to check that and ask them afterwards. It should be part of the problem ticket. I think it was a major issue from the substation that caused the power to be down.")
if (-not $ok) { Write-Output "NOT FOUND #6: A3: We'll have to check that and ask them afterwards. It should be part of the problem ticket. I think it was a major issue from the substation that caused the power to be down." } else { $rng.Text = "A3: The electricians are trying to restore power to the room with network equipment. Neerav is logging into the devices to check their status, and Praveen is checking with site users." }

$rng = $d.Content
$ok = $rng.Find.Execute("Q4: What is the status of the switches on site?")
if (-not $ok) { Write-Output "NOT FOUND #7: Q4: What is the status of the switches on site?" } else { $rng.Text = "Q4: What is the site ID for the affected location?" }

$rng = $d.Content
$ok = $rng.Find.Execute("A4: I have seen multiple switches, and all were power rebooted, but now they are reachable. Some are left, and I will update soon.")
if (-not $ok) { Write-Output "NOT FOUND #8: A4: I have seen multiple switches, and all were power rebooted, but now they are reachable. Some are left, and I will update soon." } else { $rng.Text = "A4: The site ID is XXXX." }

$rng = $d.Content
$ok = $rng.Find.Execute("Q5: Why are the ABs not registering to the controller?")
if (-not $ok) { Write-Output "NOT FOUND #9: Q5: Why are the ABs not registering to the controller?" } else { $rng.Text = "Q5: What was the root cause of the network devices going down at the Floura Site?" }

$rng = $d.Content
$ok = $rng.Find.Execute("A5: The ABs not registering to the controller are AB US-01111ab-04 on switch -sw00, AB US-01111ab-45 on switch -sw09, and AB US-01111-12 on switch -sw07.")
if (-not $ok) { Write-Output "NOT FOUND #10: A5: The ABs not registering to the controller are AB US-01111ab-04 on switch -sw00, AB US-01111ab-45 on switch -sw09, and AB US-01111-12 on switch -sw07." } else { $rng.Text = "A5: Network devices in Floura Site went down due to power issues." }

$rng = $d.Content
$ok = $rng.Find.Execute("Q6: What is the current status of the network devices in Floura Site?")
if (-not $ok) { Write-Output "NOT FOUND #11: Q6: What is the current status of the network devices in Floura Site?" } else { $rng.Text = "Q6: What actions were taken to bring the switches back online?" }

$rng = $d.Content
$ok = $rng.Find.Execute("A6: Network devices in Floura Site went down due to power issues.")
if (-not $ok) { Write-Output "NOT FOUND #12: A6: Network devices in Floura Site went down due to power issues." } else { $rng.Text = "A6: Multiple switches were power rebooted; now they are reachable. Some are left, and updates will be provided soon." }

$rng = $d.Content
$ok = $rng.Find.Execute("Q7: What was the cause of the network devices going down at the site?")
if (-not $ok) { Write-Output "NOT FOUND #13: Q7: What was the cause of the network devices going down at the site?" } else { $rng.Text = "Q7: Which ABs are not registering to the controller and on which switches are they located?" }

$rng = $d.Content
$ok = $rng.Find.Execute("A7: There was a power failure at the site, and when the power came back, there were still some power issues in the network room. The electrician was working on it, which caused the network devices to go down, isolating the site completely.")
if (-not $ok) { Write-Output "NOT FOUND #14: A7: There was a power failure at the site, and when the power came back, there were still some power issues in the network room. The electrician was working on it, which caused the network devices to go down, isolating the site completely." } else { $rng.Text = "A7: AB US-01111ab-04 is on switch -sw00; AB US-01111ab-45 is on switch -sw09; correction, AB US-01111-12 is on switch -sw07." }

$rng = $d.Content
$ok = $rng.Find.Execute("Q8: What is the current status of the network devices?")
if (-not $ok) { Write-Output "NOT FOUND #15: Q8: What is the current status of the network devices?" } else { $rng.Text = "Q8: Is the entire factory down right now?" }

$rng = $d.Content
$ok = $rng.Find.Execute("A8: The devices are up from the last 7 minutes; they have established adjacency, and most tunnels are up.")
if (-not $ok) { Write-Output "NOT FOUND #16: A8: The devices are up from the last 7 minutes; they have established adjacency, and most tunnels are up." } else { $rng.Text = "A8: It's down. The network is down. The whole plant lost power for 10 or 15 minutes. The telco room's power is still bad, affecting network gear. The factory has power, but the network is down." }

$rng = $d.Content
$ok = $rng.Find.Execute("Q9: When did the site go down and when did the devices start to ping again?")
if (-not $ok) { Write-Output "NOT FOUND #17: Q9: When did the site go down and when did the devices start to ping again?" } else { $rng.Text = "Q9: What are the applications that are getting impacted?" }

$rng = $d.Content
$ok = $rng.Find.Execute("A9: The site went down at 2:42 PM local time (18:42 GMT), and the devices started to ping again at 3:29 PM (19:29 GMT).")
if (-not $ok) { Write-Output "NOT FOUND #18: A9: The site went down at 2:42 PM local time (18:42 GMT), and the devices started to ping again at 3:29 PM (19:29 GMT)." } else { $rng.Text = "A9: All the applications." }
